$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.122973918914795
$ws.Range("B1").Value = 2.113569021224976
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.224819183349609
$ws.Range("E1").Value = 1.096145868301392
